$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values (prices) are preserved as literal text,
# matching the source data which stores them as inline strings, not numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.119.04"
$ws.Range("E2").Value = "  -0.08%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.048.89"
$ws.Range("E3").Value = "  -1.49%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.80"
$ws.Range("E5").Value = "  -2.06%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.662"
$ws.Range("E6").Value = "  -1.68%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "57.76"
$ws.Range("E7").Value = "  -1.74%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.381"
$ws.Range("E9").Value = "  -2.46%  "

$ws.Range("E10").Value = "  -2.09%  "

$ws.Range("E11").Value = "  +0.00%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.86"
$ws.Range("E12").Value = "  -0.67%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.860"
$ws.Range("E13").Value = "  +5.30%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.349.21"
$ws.Range("E14").Value = "  -1.48%  "

$ws.Range("E15").Value = "  +2.89%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.052.16"
$ws.Range("E16").Value = "  -1.35%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "17.82"
$ws.Range("E17").Value = "  +14.16%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.143.25"
$ws.Range("E18").Value = "  +0.12%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.75"
$ws.Range("E19").Value = "  +0.19%  "

$ws.Range("E20").Value = "  -3.39%  "

$ws.Range("E21").Value = "  -1.93%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "236.70"
$ws.Range("E22").Value = "  -0.98%  "

$ws.Range("E23").Value = "  +0.00%  "

$ws.Range("E24").Value = "  +1.58%  "

$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.48"
$ws.Range("E25").Value = "  +2.03%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.94"
$ws.Range("E26").Value = "  -0.03%  "

$ws.Range("E27").Value = "  -5.37%  "

$ws.Range("E28").Value = "  -1.53%  "

$ws.Range("E29").Value = "  -1.50%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.79"
$ws.Range("E30").Value = "  -0.81%  "

$ws.Range("E31").Value = "  -1.70%  "

$ws.Range("E32").Value = "  -3.18%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.46"
$ws.Range("E33").Value = "  +0.39%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0896"
$ws.Range("E34").Value = "  -2.19%  "

$ws.Range("E35").Value = "  -0.02%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.24"
$ws.Range("E36").Value = "  -2.78%  "

$ws.Range("E37").Value = "  -0.21%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.27"
$ws.Range("E38").Value = "  +16.01%  "

$ws.Range("E39").Value = "  -2.34%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.17"
$ws.Range("E40").Value = "  +16.18%  "

$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0222"
$ws.Range("E41").Value = "  -1.85%  "

$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.21"
$ws.Range("E42").Value = "  -4.40%  "

$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0958"
$ws.Range("E43").Value = "  -19.40%  "

$ws.Range("E44").Value = "  -2.72%  "

$ws.Range("E45").Value = "  -3.29%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.44"
$ws.Range("E46").Value = "  -1.05%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.274.09"
$ws.Range("E47").Value = "  -2.27%  "

$ws.Range("E48").Value = "  -3.83%  "

$ws.Range("E49").Value = "  -2.32%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.233.29"
$ws.Range("E50").Value = "  -1.43%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.72"
$ws.Range("E51").Value = "  -0.06%  "
